$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1334.0714
$ws.Range("I19").Value = 596.3333
$ws.Range("K19").Value = 596.3333
$ws.Range("M19").Value = -421.3333
$ws.Range("H51").Value = 3129.611
$ws.Range("I51").Value = 2999.9805
$ws.Range("J51").Value = 5333.3335
$ws.Range("K51").Value = 2999.9805
$ws.Range("L51").Value = 5333.3335
$ws.Range("M51").Value = -2515.9805
$ws.Range("N51").Value = -6301.3335
$ws.Range("H100").Value = 3479.3635
$ws.Range("I100").Value = 3596.625
$ws.Range("J100").Value = 3166.6667
$ws.Range("K100").Value = 3596.625
$ws.Range("L100").Value = 3166.6667
$ws.Range("M100").Value = -3055.625
$ws.Range("N100").Value = -4248.6667
$ws.Range("H103").Value = 1805.2
$ws.Range("I103").Value = 1736.1428
$ws.Range("J103").Value = 1865.625
$ws.Range("K103").Value = 5208.428400000001
$ws.Range("L103").Value = 5596.875
$ws.Range("M103").Value = -4622.428400000001
$ws.Range("N103").Value = -6768.875
$ws.Range("H106").Value = 3526.75
$ws.Range("I106").Value = 2602.1428
$ws.Range("K106").Value = 2602.1428
$ws.Range("M106").Value = -1971.1428
$ws.Range("H138").Value = 6578.625
$ws.Range("I138").Value = 3764.3333
$ws.Range("J138").Value = 6980.6665
$ws.Range("K138").Value = 11292.9999
$ws.Range("L138").Value = 20941.9995
$ws.Range("M138").Value = -6152.999899999999
$ws.Range("N138").Value = -31221.9995

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 513853.5
$ws.Range("I32").Value = 601798.9399999999
$ws.Range("K32").Value = 601798.9399999999
$ws.Range("M32").Value = -601511.9399999999
$ws.Range("H45").Value = 1947.8
$ws.Range("I45").Value = 2028.625
$ws.Range("K45").Value = 2028.625
$ws.Range("M45").Value = -1651.625
$ws.Range("H61").Value = 6244857
$ws.Range("I61").Value = 2333536.8
$ws.Range("K61").Value = 2333536.8
$ws.Range("M61").Value = -2333324.8
$ws.Range("H102").Value = 982.0526
$ws.Range("I102").Value = 964
$ws.Range("K102").Value = 964
$ws.Range("M102").Value = 658
$ws.Range("H132").Value = 2825.1667
$ws.Range("I132").Value = 1492.4688
$ws.Range("K132").Value = 4477.4064
$ws.Range("M132").Value = -1947.4064
$ws.Range("H136").Value = 6244857
$ws.Range("I136").Value = 2333536.8
$ws.Range("K136").Value = 7000610.399999999
$ws.Range("M136").Value = -6998060.399999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 41345.855
$ws.Range("I20").Value = 52399.617
$ws.Range("K20").Value = 52399.617
$ws.Range("M20").Value = -52152.617
$ws.Range("H99").Value = 12818.909
$ws.Range("I99").Value = 15312.125
$ws.Range("K99").Value = 15312.125
$ws.Range("M99").Value = -13814.125
$ws.Range("H107").Value = 2469.8333
$ws.Range("I107").Value = 2469.8333
$ws.Range("K107").Value = 2469.8333
$ws.Range("M107").Value = -549.8332999999998
$ws.Range("H134").Value = 5294254.5
$ws.Range("I134").Value = 4765118
$ws.Range("K134").Value = 14295354
$ws.Range("M134").Value = -14292819

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 7250
$ws.Range("J26").Value = 7250
$ws.Range("L26").Value = 7250
$ws.Range("N26").Value = -7824
$ws.Range("H31").Value = 783471.5
$ws.Range("I31").Value = 2948130
$ws.Range("J31").Value = 6414.59
$ws.Range("K31").Value = 2948130
$ws.Range("L31").Value = 6414.59
$ws.Range("M31").Value = -2947835
$ws.Range("N31").Value = -7004.59
$ws.Range("H34").Value = 783471.5
$ws.Range("I34").Value = 2948130
$ws.Range("J34").Value = 6414.59
$ws.Range("K34").Value = 2948130
$ws.Range("L34").Value = 6414.59
$ws.Range("M34").Value = -2947928
$ws.Range("N34").Value = -6818.59
$ws.Range("H52").Value = 79949.5
$ws.Range("J52").Value = 79949.5
$ws.Range("L52").Value = 79949.5
$ws.Range("N52").Value = -80537.5
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H132").Value = 3776.75
$ws.Range("I132").Value = 3502.2307
$ws.Range("K132").Value = 10506.6921
$ws.Range("M132").Value = -7976.6921
$ws.Range("H135").Value = 97999.5
$ws.Range("J135").Value = 97999.5
$ws.Range("L135").Value = 97999.5
$ws.Range("N135").Value = -108139.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 6107.763
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 6107.763
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H92").Value = 620.125
$ws.Range("J92").Value = 599.6667
$ws.Range("L92").Value = 1799.0001
$ws.Range("N92").Value = -4295.0001
$ws.Range("H93").Value = 2280
$ws.Range("J93").Value = 2280
$ws.Range("L93").Value = 6840
$ws.Range("N93").Value = -10584
$ws.Range("H121").Value = 82166.14999999999
$ws.Range("I121").Value = 125676.125
$ws.Range("J121").Value = 63846.156
$ws.Range("K121").Value = 377028.375
$ws.Range("L121").Value = 191538.468
$ws.Range("M121").Value = -375718.375
$ws.Range("N121").Value = -194158.468

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 6673333
$ws.Range("I18").Value = 6673333
$ws.Range("K18").Value = 6673333
$ws.Range("M18").Value = -6673040
$ws.Range("H80").Value = 7386
$ws.Range("I80").Value = 7649
$ws.Range("K80").Value = 7649
$ws.Range("M80").Value = -6651
$ws.Range("H83").Value = 7386
$ws.Range("I83").Value = 7649
$ws.Range("K83").Value = 38245
$ws.Range("M83").Value = -33253
$ws.Range("H132").Value = 10966
$ws.Range("I132").Value = 7286.7954
$ws.Range("J132").Value = 51437.25
$ws.Range("K132").Value = 21860.3862
$ws.Range("L132").Value = 154311.75
$ws.Range("M132").Value = -19330.3862
$ws.Range("N132").Value = -159371.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H30").Value = 383
$ws.Range("I30").Value = 383
$ws.Range("K30").Value = 383
$ws.Range("M30").Value = -275
$ws.Range("H68").Value = 1437.7858
$ws.Range("I68").Value = 1553.1666
$ws.Range("J68").Value = 745.5
$ws.Range("K68").Value = 1553.1666
$ws.Range("L68").Value = 745.5
$ws.Range("M68").Value = -804.1666
$ws.Range("N68").Value = -2243.5
$ws.Range("H71").Value = 1437.7858
$ws.Range("I71").Value = 1553.1666
$ws.Range("J71").Value = 745.5
$ws.Range("K71").Value = 7765.833000000001
$ws.Range("L71").Value = 3727.5
$ws.Range("M71").Value = -4021.833000000001
$ws.Range("N71").Value = -11215.5
$ws.Range("H93").Value = 7832.6665
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 7832.6665
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -10328.6665
$ws.Range("H127").Value = 154357.5
$ws.Range("J127").Value = 154357.5
$ws.Range("L127").Value = 154357.5
$ws.Range("N127").Value = -164277.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5644
$ws.Range("J62").Value = 4724.25
$ws.Range("L62").Value = 4724.25
$ws.Range("N62").Value = -5972.25
$ws.Range("H65").Value = 5644
$ws.Range("J65").Value = 4724.25
$ws.Range("L65").Value = 23621.25
$ws.Range("N65").Value = -29861.25
$ws.Range("H122").Value = 40443.242
$ws.Range("I122").Value = 1207.2632
$ws.Range("J122").Value = 114991.6
$ws.Range("K122").Value = 3621.7896
$ws.Range("L122").Value = 344974.8
$ws.Range("M122").Value = -1171.7896
$ws.Range("N122").Value = -349874.8
$ws.Range("H132").Value = 5379513.5
$ws.Range("I132").Value = 5955518.5
$ws.Range("J132").Value = 3466.3333
$ws.Range("K132").Value = 17866555.5
$ws.Range("L132").Value = 10398.9999
$ws.Range("M132").Value = -17864025.5
$ws.Range("N132").Value = -15458.9999
$ws.Range("H136").Value = 2979874
$ws.Range("J136").Value = 4814154
$ws.Range("L136").Value = 14442462
$ws.Range("N136").Value = -14447562
